# [ANV] updating decay chain spreadsheet
# Adds a new "Cu Density" worksheet (mirroring the existing "HDPE Density"
# sheet) right before the existing "Cu Target Fractions " sheet.

$wb = $excel.ActiveWorkbook

$hdpe = $wb.Worksheets.Item("HDPE Density")
$cuTargetFractions = $wb.Worksheets.Item("Cu Target Fractions ")

# Insert the new sheet immediately before "Cu Target Fractions "
$newSheet = $wb.Worksheets.Add($cuTargetFractions)
$newSheet.Name = "Cu Density"

# --- column widths (mirror HDPE Density) ---
for ($i = 1; $i -le 9; $i++) {
    $newSheet.Columns.Item($i).ColumnWidth = $hdpe.Columns.Item($i).ColumnWidth
}

# --- copy row formatting (fills/borders) from the analogous HDPE rows ---
$hdpe.Range("A1:I1").Copy()
$newSheet.Range("A1:I1").PasteSpecial(-4122)   # xlPasteFormats

$hdpe.Range("A4:I4").Copy()
$newSheet.Range("A3:I3").PasteSpecial(-4122)   # xlPasteFormats (Total row)

$hdpe.Range("A12:I12").Copy()
$newSheet.Range("A11:I11").PasteSpecial(-4122) # xlPasteFormats (sub-table header)

$excel.CutCopyMode = 0

# --- header row ---
$newSheet.Range("A1").Value = "Element"
$newSheet.Range("B1").Value = "Z (atomic number)"
$newSheet.Range("C1").Value = "Mass Fraction (%)"
$newSheet.Range("D1").Value = "overall density (kg/m^3)"
$newSheet.Range("E1").Value = "overall density (g/cm^3)"
$newSheet.Range("F1").Value = "Elemental Mass Density (g/cm^3)"
$newSheet.Range("G1").Value = "Molar Mass (amu)"
$newSheet.Range("H1").Value = "Number Density (#/cm^3)"
$newSheet.Range("I1").Value = "fraction of atoms"

# --- Cu data row ---
$newSheet.Range("A2").Value = "Cu"
$newSheet.Range("B2").Value = 6
$newSheet.Range("C2").Value = 1
$newSheet.Range("D2").Value = 8960
$newSheet.Range("E2").Formula = "=D2*(1000)*(1/100000)"
$newSheet.Range("F2").Formula = "=`$E`$2*(C2/100)"
$newSheet.Range("G2").Value = 63.546
$newSheet.Range("H2").Formula = "=(F2/G2)*6.0221408E+23"
$newSheet.Range("I2").Formula = "=H2/`$H`$3"

# --- Total row ---
$newSheet.Range("A3").Value = "Total"
$newSheet.Range("C3").Formula = "=SUM(C2:C2)"
$newSheet.Range("F3").Formula = "=SUM(F2:F2)"
$newSheet.Range("H3").Formula = "=SUM(H2:H2)"
$newSheet.Range("I3").Formula = "=SUM(I2:I2)"

# --- reference link ---
$newSheet.Range("A6").Value = "https://en.wikipedia.org/wiki/Copper"

# --- small summary sub-table ---
$newSheet.Range("A11").Value = "Z (atomic number)"
$newSheet.Range("B11").Value = "fraction of atoms"
$newSheet.Range("A12").Value = 29
$newSheet.Range("B12").Value = 1

# --- view settings for the new sheet ---
$win = $excel.ActiveWindow
$win.Zoom = 101
$newSheet.Range("B25").Select()

$wb.Save()
